$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_14")

# --- Remove the yellow highlight fill from the data table (A2:K21) ---
$ws.Range("A2:K21").Interior.ColorIndex = -4142

# --- Add the new "Total check" row (row 22) ---
$ws.Range("A22").Value = "Total check"
$ws.Range("A22").Font.Color = 255

$ws.Range("B22").Formula = "=SUM(B2:B19)-B20"
$ws.Range("C22:K22").Formula = "=SUM(C2:C19)-C20"

$ws.Range("B22:K22").Font.Color = 255
$ws.Range("B22:K22").NumberFormat = "#,##0"
